$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "33÷3=11, 0"
$t.Cell(1, 2).Range.Text = "12÷6=2, 0"
$t.Cell(1, 3).Range.Text = "53÷6=8, 5"
$t.Cell(1, 4).Range.Text = "66÷6=11, 0"
$t.Cell(1, 5).Range.Text = "36÷7=5, 1"
$t.Cell(5, 1).Range.Text = "31÷7=4, 3"
$t.Cell(5, 2).Range.Text = "67÷7=9, 4"
$t.Cell(5, 3).Range.Text = "47÷4=11, 3"
$t.Cell(5, 4).Range.Text = "49÷3=16, 1"
$t.Cell(5, 5).Range.Text = "74÷3=24, 2"
$t.Cell(9, 1).Range.Text = "25÷2=12, 1"
$t.Cell(9, 2).Range.Text = "35÷2=17, 1"
$t.Cell(9, 3).Range.Text = "68÷7=9, 5"
$t.Cell(9, 4).Range.Text = "59÷3=19, 2"
$t.Cell(9, 5).Range.Text = "59÷4=14, 3"
$t.Cell(13, 1).Range.Text = "53÷8=6, 5"
$t.Cell(13, 2).Range.Text = "81÷8=10, 1"
$t.Cell(13, 3).Range.Text = "43÷7=6, 1"
$t.Cell(13, 4).Range.Text = "90÷7=12, 6"
$t.Cell(13, 5).Range.Text = "45÷6=7, 3"
$t.Cell(17, 1).Range.Text = "61÷2=30, 1"
$t.Cell(17, 2).Range.Text = "77÷3=25, 2"
$t.Cell(17, 3).Range.Text = "57÷8=7, 1"
$t.Cell(17, 4).Range.Text = "79÷9=8, 7"
$t.Cell(17, 5).Range.Text = "83÷4=20, 3"
